$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Write notation table content (rows 1-21) in row-major order ---
$ws1.Range("A1").Value = 'Current'
$ws1.Range("B1").Value = 'Item'
$ws1.Range("C1").Value = 'Final'
$ws1.Range("D1").Value = 'Reason'

$ws1.Range("A2").Value = 'D'
$ws1.Range("B2").Value = 'data matrix '
$ws1.Range("C2").Value = 'D'
$ws1.Range("D2").Value = 'D for data, matrix so capital'

$ws1.Range("A3").Value = 'd'
$ws1.Range("B3").Value = 'data vector'
$ws1.Range("C3").Value = 'd'
$ws1.Range("D3").Value = 'd for data, vector so lower case'

$ws1.Range("A4").Value = 'n'
$ws1.Range("B4").Value = 'number of samples'
$ws1.Range("C4").Value = 'n'
$ws1.Range("D4").Value = 'n for number, widespread use in science'

$ws1.Range("A5").Value = 'p'
$ws1.Range("B5").Value = 'p - number of variables - change to m?'
$ws1.Range("C5").Value = 'm'
$ws1.Range("D5").Value = 'm for nuMber (and adjacent to n, mxn is a common convention for matrix dimensions)'

$ws1.Range("A6").Value = 'k'
$ws1.Range("B6").Value = 'sample index'
$ws1.Range("C6").Value = 'o'
$ws1.Range("D6").Value = 'o for observation'

$ws1.Range("A7").Value = 'j'
$ws1.Range("B7").Value = 'variable index - j is currentlly used in two ways'
$ws1.Range("C7").Value = 'v'
$ws1.Range("D7").Value = 'v for variable'

$ws1.Range("A8").Value = 'i'
$ws1.Range("B8").Value = 'PC index'
$ws1.Range("C8").Value = 'p'
$ws1.Range("D8").Value = 'p for Pc, scalar'

$ws1.Range("A9").Value = 'd'
$ws1.Range("B9").Value = 'Number of PCs/dimensionality'
$ws1.Range("C9").Value = 'f'
$ws1.Range("D9").Value = 'f for number of factors'

$ws1.Range("A10").Value = 'S'
$ws1.Range("B10").Value = 'Score matrix'
$ws1.Range("C10").Value = 'S'
$ws1.Range("D10").Value = 'S for score, matrix'

$ws1.Range("A11").Value = 's'
$ws1.Range("B11").Value = 'score vector'
$ws1.Range("C11").Value = 's'
$ws1.Range("D11").Value = 's for score, vector or scalar'

$ws1.Range("A12").Value = 'L'
$ws1.Range("B12").Value = 'loading matrix (ensure denoted with transpose where required)'
$ws1.Range("C12").Value = 'L'
$ws1.Range("D12").Value = 'L for loading, matrix'

$ws1.Range("A13").Value = 'l'
$ws1.Range("B13").Value = 'loading vector'
$ws1.Range("C13").Value = 'l'
$ws1.Range("D13").Value = 'l for loading, vector'

$ws1.Range("A14").Value = '\Top'
$ws1.Range("B14").Value = 'character used for transpose'
$ws1.Range("C14").Value = '^\top'
$ws1.Range("D14").Value = 'widespread convention'

$ws1.Range("A15").Value = '^\dagger'
$ws1.Range("B15").Value = 'character used for pseudoinverse'
$ws1.Range("C15").Value = '^\dagger'
$ws1.Range("D15").Value = 'widespread convention'

$ws1.Range("A16").Value = '*'
$ws1.Range("B16").Value = 'multiplication'
$ws1.Range("C16").Value = '\times'
$ws1.Range("D16").Value = 'conventional mathemetical notation'

$ws1.Range("B17").Value = 'Conditional statement'
$ws1.Range("C17").Value = '[]'
$ws1.Range("D17").Value = 'Iverson Brackets'

$ws1.Range("A18").Value = 'j'
$ws1.Range("B18").Value = 'iteration in NIPALS'
$ws1.Range("C18").Value = 'i'
$ws1.Range("D18").Value = 'i for iteration'

$ws1.Range("B19").Value = 'elements conditional on positive score'
$ws1.Range("C19").Value = '_{p+}'
$ws1.Range("D19").Value = 'loading or score for pth PC from positive score spectra'

$ws1.Range("B20").Value = 'elements conditional on negative score'
$ws1.Range("C20").Value = '_{p-}'
$ws1.Range("D20").Value = 'loading or score for pth PC from negative score spectra'

$ws1.Range("B21").Value = 'covariance matrix'
$ws1.Range("C21").Value = 'C'

# --- Apply theme font color (style index 5) to column A, rows 1-21 ---
$ws1.Range("A1:A21").Font.ThemeColor = 4

# --- Selections: set sheet1 selection first, then sheet2 last so sheet2 stays active tab ---
$ws1.Range("B28").Select()
$ws2.Range("D5").Select()

Write-Output "edit complete"